# Weekly fruit/vegetable price update.
# A new week's price record is inserted above the current row 22, shifting the
# existing rows 22-24 down to 23-25, and the new row 22 is populated with the
# latest week's data (2021-08-12 / serial 44420).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 22 (pushes old rows 22..24 to 23..25).
$ws.Rows.Item(22).Insert()

# Populate the new row 22 with the new weekly record.
$ws.Range("A22").Value = 10
$ws.Range("B22").Value = "Vega Modelo de Temuco"
$ws.Range("C22").Value = "La Araucanía"
$ws.Range("D22").Value = 44420
$ws.Range("E22").Value = 9
$ws.Range("F22").Value = 100112010
$ws.Range("G22").Value = "Achicoria"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 45
$ws.Range("K22").Value = 8000
$ws.Range("L22").Value = 8000
$ws.Range("M22").Value = 8000
$ws.Range("N22").Value = "`$/caja 16 unidades"
$ws.Range("O22").Value = "Región Metropolitana"
$ws.Range("P22").Value = 500
$ws.Range("Q22").Value = 16
$ws.Range("R22").Value = "Hortaliza"
